# Fix wrong url redirect issue:
#  - insert a new header row ("Name" / "Nicknames") above the data
#  - remove the (broken) hyperlinks that were attached to the old A-column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 1 and shift all existing data down.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Nicknames"

# Drop the (mis-pointed) hyperlinks entirely.
$ws.Hyperlinks.Delete()

# Leave the selection on the first data row, like the author's saved file.
$ws.Range("A2").Select()
